{"js": "const tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst t = tables.items[0];\nt.load(\"values\");\nawait context.sync();\nconst rowCount = t.values.length;\nconst columnCount = t.values[0].length;\n\n// New expressions, row-major, matching the table's existing 20x5 grid.\nconst newValues = [\n  [\"74+11=\", \"95-67=\", \"80-43=\", \"98-70=\", \"44+25=\"],\n  [\"60-33=\", \"8+16=\", \"68-61=\", \"51-22=\", \"42-11=\"],\n  [\"26+52=\", \"4+56=\", \"12+29=\", \"46+16=\", \"62-13=\"],\n  [\"74-63=\", \"9-8=\", \"66-50=\", \"80-46=\", \"3+62=\"],\n  [\"68+23=\", \"85-18=\", \"4+74=\", \"3+43=\", \"66+11=\"],\n  [\"68-52=\", \"52-45=\", \"5+92=\", \"75-20=\", \"92-83=\"],\n  [\"79-21=\", \"42+23=\", \"76+0=\", \"80-68=\", \"52+5=\"],\n  [\"60+37=\", \"68+5=\", \"8+51=\", \"82-52=\", \"44+2=\"],\n  [\"32+44=\", \"33+26=\", \"99-24=\", \"63-17=\", \"81-43=\"],\n  [\"27+59=\", \"6+25=\", \"89-59=\", \"95-55=\", \"76-0=\"],\n  [\"44-17=\", \"93-0=\", \"56-41=\", \"83-33=\", \"90-26=\"],\n  [\"78-76=\", \"75+13=\", \"15+59=\", \"29-9=\", \"26+39=\"],\n  [\"68-6=\", \"4+9=\", \"3+30=\", \"6-4=\", \"36+48=\"],\n  [\"38+32=\", \"43-42=\", \"1+6=\", \"42-2=\", \"30+42=\"],\n  [\"90-77=\", \"84-47=\", \"68+0=\", \"77-11=\", \"38+9=\"],\n  [\"94-65=\", \"72+0=\", \"70-11=\", \"62-6=\", \"26+20=\"],\n  [\"11+60=\", \"28-5=\", \"75-44=\", \"38-4=\", \"67-20=\"],\n  [\"31+6=\", \"3+9=\", \"28+47=\", \"43+5=\", \"86-67=\"],\n  [\"27-11=\", \"40+26=\", \"9+23=\", \"94-89=\", \"44+55=\"],\n  [\"26+6=\", \"27+48=\", \"26+64=\", \"22+5=\", \"93+1=\"]\n];\n\nconst rows = Math.min(rowCount, newValues.length);\nfor (let r = 0; r < rows; r++) {\n  const cols = Math.min(columnCount, newValues[r].length);\n  for (let c = 0; c < cols; c++) {\n    const cell = t.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$values = @(\n    '74+11=',\n    '95-67=',\n    '80-43=',\n    '98-70=',\n    '44+25=',\n    '60-33=',\n    '8+16=',\n    '68-61=',\n    '51-22=',\n    '42-11=',\n    '26+52=',\n    '4+56=',\n    '12+29=',\n    '46+16=',\n    '62-13=',\n    '74-63=',\n    '9-8=',\n    '66-50=',\n    '80-46=',\n    '3+62=',\n    '68+23=',\n    '85-18=',\n    '4+74=',\n    '3+43=',\n    '66+11=',\n    '68-52=',\n    '52-45=',\n    '5+92=',\n    '75-20=',\n    '92-83=',\n    '79-21=',\n    '42+23=',\n    '76+0=',\n    '80-68=',\n    '52+5=',\n    '60+37=',\n    '68+5=',\n    '8+51=',\n    '82-52=',\n    '44+2=',\n    '32+44=',\n    '33+26=',\n    '99-24=',\n    '63-17=',\n    '81-43=',\n    '27+59=',\n    '6+25=',\n    '89-59=',\n    '95-55=',\n    '76-0=',\n    '44-17=',\n    '93-0=',\n    '56-41=',\n    '83-33=',\n    '90-26=',\n    '78-76=',\n    '75+13=',\n    '15+59=',\n    '29-9=',\n    '26+39=',\n    '68-6=',\n    '4+9=',\n    '3+30=',\n    '6-4=',\n    '36+48=',\n    '38+32=',\n    '43-42=',\n    '1+6=',\n    '42-2=',\n    '30+42=',\n    '90-77=',\n    '84-47=',\n    '68+0=',\n    '77-11=',\n    '38+9=',\n    '94-65=',\n    '72+0=',\n    '70-11=',\n    '62-6=',\n    '26+20=',\n    '11+60=',\n    '28-5=',\n    '75-44=',\n    '38-4=',\n    '67-20=',\n    '31+6=',\n    '3+9=',\n    '28+47=',\n    '43+5=',\n    '86-67=',\n    '27-11=',\n    '40+26=',\n    '9+23=',\n    '94-89=',\n    '44+55=',\n    '26+6=',\n    '27+48=',\n    '26+64=',\n    '22+5=',\n    '93+1='\n)\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        if ($idx -ge $values.Length) { break }\n        $cell = $t.Cell($r, $c)\n        $cell.Range.Text = $values[$idx]\n        $idx++\n    }\n}"}
